$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the containing folder path (e.g. "Brent_Large/") from each
# fastqFileName value in column L (rows 2-57), leaving just the bare
# filename.
for ($row = 2; $row -le 57; $row++) {
    $cell = $ws.Cells.Item($row, 12)
    $current = [string]$cell.Text
    $slashIndex = $current.IndexOf("/")
    if ($slashIndex -ge 0) {
        $newValue = $current.Substring($slashIndex + 1)
        $cell.Value = $newValue
    }
}
